$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'247.15"
$ws.Range("E2").Value = "'0.98%"
$ws.Range("D3").Value = "'26.25"
$ws.Range("E3").Value = "'4.54%"
$ws.Range("D4").Value = "'5.086"
$ws.Range("E4").Value = "'1.09%"
$ws.Range("D5").Value = "'0.05601"
$ws.Range("E5").Value = "'-0.19%"
$ws.Range("D6").Value = "'6.476"
$ws.Range("E6").Value = "'-1.47%"
$ws.Range("E7").Value = "'-0.10%"
$ws.Range("D8").Value = "'0.8461"
$ws.Range("E8").Value = "'1.19%"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "'0.1338"
$ws.Range("E9").Value = "'-0.03%"
$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D10").Value = "'0.06996"
$ws.Range("E10").Value = "'0.57%"
$ws.Range("B11").Value = "BitrueCoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D11").Value = "'0.02810"
$ws.Range("E11").Value = "'-1.02%"
$ws.Range("B12").Value = "BitMartToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D12").Value = "'0.09385"
$ws.Range("E12").Value = "'-0.19%"
$ws.Range("B13").Value = "BitForexToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D13").Value = "'0.001520"
$ws.Range("E13").Value = "'0.40%"
$ws.Range("B14").Value = "One"
$ws.Range("C14").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D14").Value = "'0.0005980"
$ws.Range("E14").Value = "'-93.84%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.006151"
$ws.Range("E15").Value = "'-0.63%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.605"
$ws.Range("E16").Value = "'3.00%"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'3.020"
$ws.Range("E17").Value = "'0.33%"
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").Value = "'2.055"
$ws.Range("E18").Value = "'-1.72%"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").Value = "'0.3113"
$ws.Range("E19").Value = "'-1.47%"
$ws.Range("D20").Value = "'0.03190"
$ws.Range("E20").Value = "'-1.99%"
$ws.Range("E21").Value = "'-1.30%"
$ws.Range("D22").Value = "'3.743"
$ws.Range("E22").Value = "'-0.01%"
$ws.Range("D23").Value = "'0.04656"
$ws.Range("E23").Value = "'-0.62%"
$ws.Range("D24").Value = "'0.1374"
$ws.Range("E24").Value = "'0.37%"
$ws.Range("D25").Value = "'0.001243"
$ws.Range("E25").Value = "'0.05%"
$ws.Range("D27").Value = "'0.00009601"
$ws.Range("E27").Value = "'-0.95%"
$ws.Range("E28").Value = "'-0.02%"
$ws.Range("D40").Value = "'0.03663"
$ws.Range("E40").Value = "'0.04%"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.006150"
$ws.Range("E41").Value = "'-1.04%"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1057"
$ws.Range("E42").Value = "'0.29%"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.002500"
$ws.Range("E43").Value = "'-8.52%"
$ws.Range("D44").Value = "'0.008269"
$ws.Range("E44").Value = "'1.08%"
$ws.Range("D45").Value = "'0.00005393"
$ws.Range("E45").Value = "'1.91%"
$ws.Range("E46").Value = "'0.08%"
$ws.Range("E47").Value = "'-35.79%"
$ws.Range("D48").Value = "'0.002431"
$ws.Range("E48").Value = "'20.54%"
$ws.Range("E49").Value = "'0.08%"
$ws.Range("E50").Value = "'0.08%"
